# Updated cryptos list on Sun Jul  7 21:39:05 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns with new scraped
# figures, and fixes the WrappedEther/Polkadot row ordering (rows 17-18
# had swapped data).
#
# Note: several "Price" values look numeric (e.g. "133.83") but must stay
# plain text, matching the source sheet's inlineStr cells. Assigning such
# a literal directly would make Excel auto-convert it to a number, so
# those are written with a leading apostrophe (forces text entry) and the
# cell style is then reset to "Normal" so no stray quote-prefix formatting
# is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '56.275.57'
$ws.Range('E2').Value = '  -2.84%  '
$ws.Range('D3').Value = '2.943.45'
$ws.Range('E3').Value = '  -4.03%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'492.59"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -6.48%  '
$ws.Range('D6').Value = "'133.83"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -6.83%  '
$ws.Range('E7').Value = '  +0.09%  '
$ws.Range('E8').Value = '  -5.88%  '
$ws.Range('E9').Value = '  -6.87%  '
$ws.Range('E10').Value = '  -7.16%  '
$ws.Range('D12').Value = '3.467.98'
$ws.Range('E12').Value = '  -3.49%  '
$ws.Range('E13').Value = '  -3.20%  '
$ws.Range('D14').Value = "'25.64"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -6.50%  '
$ws.Range('E15').Value = '  -8.93%  '
$ws.Range('D16').Value = '56.416.74'
$ws.Range('E16').Value = '  -2.62%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D17').Value = "'5.95"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -4.38%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '2.950.11'
$ws.Range('E18').Value = '  -3.88%  '
$ws.Range('E19').Value = '  -6.57%  '
$ws.Range('E20').Value = '  -6.26%  '
$ws.Range('D21').Value = "'315.64"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -7.64%  '
$ws.Range('D22').Value = "'1.00"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.03%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('E24').Value = '  -4.36%  '
$ws.Range('D25').Value = "'62.13"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.35%  '
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('E27').Value = '  -6.09%  '
$ws.Range('D28').Value = '0.0₃0859'
$ws.Range('E28').Value = '  -12.30%  '
$ws.Range('D29').Value = "'6.42"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -8.09%  '
$ws.Range('E30').Value = '  -7.20%  '
$ws.Range('E31').Value = '  -6.73%  '
$ws.Range('D32').Value = "'19.80"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.26%  '
$ws.Range('D34').Value = "'152.96"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.53%  '
$ws.Range('D35').Value = "'4.44"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -7.42%  '
$ws.Range('D36').Value = "'5.64"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.16%  '
$ws.Range('E37').Value = '  -9.67%  '
$ws.Range('D38').Value = "'23.55"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -10.75%  '
$ws.Range('E39').Value = '  -8.17%  '
$ws.Range('D40').Value = "'37.51"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.97%  '
$ws.Range('D41').Value = '2.976.08'
$ws.Range('E41').Value = '  -4.10%  '
$ws.Range('E42').Value = '  +0.05%  '
$ws.Range('E43').Value = '  -4.62%  '
$ws.Range('D44').Value = "'3.64"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.37%  '
$ws.Range('D45').Value = '2.138.56'
$ws.Range('E45').Value = '  -8.63%  '
$ws.Range('E46').Value = '  -9.45%  '
$ws.Range('D47').Value = "'5.83"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.65%  '
$ws.Range('E48').Value = '  -11.53%  '
$ws.Range('E49').Value = '  -6.82%  '
$ws.Range('D50').Value = "'18.73"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -7.33%  '
$ws.Range('E51').Value = '  -5.34%  '
